$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 526, shifting existing rows 526:622 down to 527:623.
$ws.Rows(526).Insert()

# Populate the newly inserted row 526 with its data. The categorical /
# template columns (A,B,C,E,F,G,H,I,J,K,L,Q,T) are identical to the rest of
# this "Mango - Vega Central Mapocho de Santiago" block; only the date,
# volume/price and origin columns (D,M,N,O,P,R,S) hold new data.
$ws.Cells.Item(526,1).Value  = 9
$ws.Cells.Item(526,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(526,3).Value  = "Metropolitana"
$ws.Cells.Item(526,4).Value  = 44995
$ws.Cells.Item(526,5).Value  = 13
$ws.Cells.Item(526,6).Value  = "Fruta"
$ws.Cells.Item(526,7).Value  = 100108
$ws.Cells.Item(526,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(526,9).Value  = 100108002
$ws.Cells.Item(526,10).Value = "Mango"
$ws.Cells.Item(526,11).Value = "Sin especificar"
$ws.Cells.Item(526,12).Value = "Primera"
$ws.Cells.Item(526,13).Value = 580
$ws.Cells.Item(526,14).Value = 6000
$ws.Cells.Item(526,15).Value = 6500
$ws.Cells.Item(526,16).Value = 6241
$ws.Cells.Item(526,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(526,18).Value = "Perú"
$ws.Cells.Item(526,19).Value = 1560
$ws.Cells.Item(526,20).Value = 4
